# Tested different feat. types and sim. metrics on the instance-based VCP notebook
#
# The "KNN(Instance-based)" sheet is the active sheet/tab. It holds two
# stacked tables (RMSE and MICROSEGUNDOS), each with rows per feature
# type. The HIST-EUCLID / HIST-COS rows (6-7 and 14-15) were previously
# blank except for their row labels; this run fills in the measured
# results for those two feature types across all k columns (B:F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- RMSE results table ---------------------------------------------
# Row 6: HIST-EUCLID
$ws.Range("B6").Value = 6.3
$ws.Range("C6").Value = 4.99
$ws.Range("D6").Value = 4.6399999999999997
$ws.Range("E6").Value = 4.55
$ws.Range("F6").Value = 4.41

# Row 7: HIST-COS
$ws.Range("B7").Value = 6.39
$ws.Range("C7").Value = 5.04
$ws.Range("D7").Value = 4.7
$ws.Range("E7").Value = 4.4800000000000004
$ws.Range("F7").Value = 4.3899999999999997

# --- MICROSEGUNDOS (timing) results table ----------------------------
# Row 14: HIST-EUCLID
$ws.Range("B14").Value = 3783
$ws.Range("C14").Value = 3957
$ws.Range("D14").Value = 3938
$ws.Range("E14").Value = 4337
$ws.Range("F14").Value = 5182

# Row 15: HIST-COS
$ws.Range("B15").Value = 6934
$ws.Range("C15").Value = 6944
$ws.Range("D15").Value = 7056
$ws.Range("E15").Value = 7368
$ws.Range("F15").Value = 8232

# Leave the selection where the author ended up after entering the data.
$ws.Range("D19").Select()
